# Apply the cryptos list refresh (prices + 1h change %) described by the diff.
#
# All data cells in this sheet are stored as literal text (t="inlineStr"),
# even when their content looks like a plain number (e.g. "584.30"). Assigning
# such a string straight to Range.Value would make Excel auto-convert it to a
# numeric value, same as typing it into a General-formatted cell. To keep those
# specific cells as text (matching the source data / diff) we temporarily mark
# them as Text-formatted before the write, then restore the default "Normal"
# style so no visible formatting change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.537.62"
$ws.Range("E2").Value = "  +5.20%  "

# Row 3
$ws.Range("D3").Value = "3.092.12"
$ws.Range("E3").Value = "  +2.95%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.03%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "3.087.80"
$ws.Range("E8").Value = "  +3.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.52%  "

# Row 15
$ws.Range("E15").Value = "  +0.70%  "

# Row 16
$ws.Range("D16").Value = "3.599.90"
$ws.Range("E16").Value = "  +2.87%  "

# Row 17
$ws.Range("E17").Value = "  +0.18%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.438.90"
$ws.Range("E18").Value = "  +5.11%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.082.80"
$ws.Range("E19").Value = "  +2.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.26%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("E31").Value = "  +12.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.92%  "

# Row 34
$ws.Range("E34").Value = "  +4.58%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("E35").Value = "  +4.93%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "425.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.08%  "

# Row 42
$ws.Range("D42").Value = "2.925.06"
$ws.Range("E42").Value = "  +6.21%  "

# Row 43
$ws.Range("E43").Value = "  +4.51%  "

# Row 44
$ws.Range("E44").Value = "  +10.36%  "

# Row 45
$ws.Range("E45").Value = "  +0.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.04%  "

# Row 47
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.01%  "
